$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 45

$ws.Cells.Item($row, 1).Value = 44
$ws.Cells.Item($row, 2).Value = "lebanon"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45240.625
$ws.Cells.Item($row, 6).Value = "Safa"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Bourj FC"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 2.38
$ws.Cells.Item($row, 11).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 12).Value = 2.38
$ws.Cells.Item($row, 13).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 14).Value = 3.25
$ws.Cells.Item($row, 15).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 16).Value = 3.25
$ws.Cells.Item($row, 17).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 18).Value = 2.92
$ws.Cells.Item($row, 19).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 20).Value = 2.92
$ws.Cells.Item($row, 21).Value = "10/11/2023 11:52"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/lebanon/premier-league/safa-bourj/SrCeANgo/"

# Copy style from row 44 (A44 -> A45, E44 -> E45) to preserve formatting
$ws.Cells.Item(44, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Cells.Item(44, 5).Copy() | Out-Null
$ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null # xlPasteFormats
